$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("供应链放还款202312")

$ws.Range("H2").Value = 4865.75
$ws.Range("H3").Value = 1013.94
$ws.Range("H4").Value = 5975.78
$ws.Range("H5").Value = 174.82
$ws.Range("H6").Value = 78.56
$ws.Range("H7").Value = 18837.17
$ws.Range("H8").Value = 75743.83
$ws.Range("H9").Value = 91.54
$ws.Range("H10").Value = 7.12
$ws.Range("H11").Value = 8.22
$ws.Range("H12").Value = 16626.24
$ws.Range("H13").Value = 18782.94
$ws.Range("H14").Value = 20342.45
$ws.Range("H15").Value = 227509.14
$ws.Range("H16").Value = 1453.31
$ws.Range("H17").Value = 281.57
$ws.Range("H18").Value = 391792.38
